$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the total for column D (sum of captures) in row 7, matching the
# existing SUM(C2:C6) pattern already present in C7.
$ws.Range("D7").Formula() = "=SUM(D2:D6)"

# New header for column G.
$ws.Range("G1").Value() = "Cap. Percent (%)"

# New column G: percent of total captures, one formula per row (relative
# row reference, absolute reference to the D7 total) -- Excel will group
# the fill-down formulas G3:G6 into a shared formula automatically.
$ws.Range("G2").Formula() = '=(D2/$D$7)*100'
$ws.Range("G3:G6").Formula() = '=(D3/$D$7)*100'

# Give the new column a width similar to the others.
$ws.Columns.Item(7).ColumnWidth() = 20.830729166666668

# Recalculate so the cached formula values are stored in the file.
$excel.Calculate()

# Match the author's final cell selection.
$ws.Range("G13").Select()
